$d = $word.ActiveDocument

# 1. Color the "Use the VS solution ..." paragraph (Exceptions section) red.
$d.Paragraphs(7).Range.Font.Color = 255

# 2. Color the "Add exception handling ..." bullet paragraph (Exceptions section) red.
$d.Paragraphs(8).Range.Font.Color = 255

# 3. Move the "_GoBack" bookmark from the end of the "Assertions" section
#    paragraph to the very last (empty) paragraph of the document. Adding a
#    bookmark with the same name removes any pre-existing bookmark of that
#    name, so this both deletes the old one and creates the new one.
$lastParagraph = $d.Paragraphs.Last
$d.Bookmarks.Add("_GoBack", $lastParagraph.Range)
